# The project is ready to present
# Clean up the product catalogue: normalise IDs to sequential numbers,
# normalise stock_type wording ("kg" -> "Unpackaged", "pack" -> "Packaged"),
# strip the unit suffix out of stock_amount so it is numeric, remove the
# left-over scratch rows, add the missing Cherry entry, and rename the
# category tabs to their fuller names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Fruits
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "apple"
$ws1.Range("C2").Value = "-"
$ws1.Range("D2").Value = "Unpackaged"
$ws1.Range("E2").Value = 1.2
$ws1.Range("F2").Value = 90

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Mandarin"
$ws1.Range("C3").Value = "-"
$ws1.Range("D3").Value = "Unpackaged"
$ws1.Range("E3").Value = 2.3
$ws1.Range("F3").Value = 100

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "Portugal"
$ws1.Range("C4").Value = "-"
$ws1.Range("D4").Value = "Unpackaged"
$ws1.Range("E4").Value = 1.85
$ws1.Range("F4").Value = 85

$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "Banana"
$ws1.Range("C5").Value = "-"
$ws1.Range("D5").Value = "Unpackaged"
$ws1.Range("E5").Value = 5.85
$ws1.Range("F5").Value = 50

$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = "Strawberry"
$ws1.Range("C6").Value = "-"
$ws1.Range("D6").Value = "Unpackaged"
$ws1.Range("E6").Value = 9.85
$ws1.Range("F6").Value = 50

# Drop the two leftover scratch rows and turn row 7 into a real product
$ws1.Rows.Item(9).Delete()
$ws1.Rows.Item(8).Delete()

$ws1.Range("A7").Value = 7
$ws1.Range("B7").Value = "Cherry"
$ws1.Range("C7").Value = "-"
$ws1.Range("D7").Value = "Unpackaged"
$ws1.Range("E7").Value = 10
$ws1.Range("F7").Value = 30

$ws1.Activate()
$ws1.Rows.Item(7).Select()

# ---------------------------------------------------------------------
# Sheet 2: Vegetables
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = 1
$ws2.Range("D2").Value = "Unpackaged"
$ws2.Range("F2").Value = 90

$ws2.Range("A3").Value = 2
$ws2.Range("D3").Value = "Unpackaged"
$ws2.Range("F3").Value = 100

$ws2.Range("A4").Value = 3
$ws2.Range("D4").Value = "Unpackaged"
$ws2.Range("F4").Value = 50

$ws2.Range("A5").Value = 4
$ws2.Range("D5").Value = "Unpackaged"
$ws2.Range("F5").Value = 50

$ws2.Range("A6").Value = 5
$ws2.Range("D6").Value = "Unpackaged"
$ws2.Range("F6").Value = 80

# Remove the leftover scratch row
$ws2.Rows.Item(7).Delete()

$ws2.Activate()
$ws2.Rows.Item(7).Select()

# ---------------------------------------------------------------------
# Sheet 3: Dairy -> Dairies and Cereal
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Dairies and Cereal"

$ws3.Range("A2").Value = 1
$ws3.Range("D2").Value = "Packaged"

$ws3.Range("A3").Value = 2
$ws3.Range("D3").Value = "Unpackaged"
$ws3.Range("F3").Value = 50

$ws3.Range("A4").Value = 3
$ws3.Range("D4").Value = "Unpackaged"
$ws3.Range("F4").Value = 50

$ws3.Range("A5").Value = 4
$ws3.Range("D5").Value = "Packaged"

$ws3.Range("A6").Value = 5
$ws3.Range("D6").Value = "Packaged"

$ws3.Activate()
$ws3.Range("D4").Select()

# ---------------------------------------------------------------------
# Sheet 4: Meat -> Meat Products
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "Meat Products"

$ws4.Range("A2").Value = 1
$ws4.Range("D2").Value = "Packaged"

$ws4.Range("A3").Value = 2
$ws4.Range("D3").Value = "Packaged"
$ws4.Range("E3").Value = 40.3

$ws4.Range("A4").Value = 3
$ws4.Range("D4").Value = "Unpackaged"
$ws4.Range("F4").Value = 50

$ws4.Range("A5").Value = 4
$ws4.Range("D5").Value = "Packaged"

$ws4.Range("A6").Value = 5
$ws4.Range("D6").Value = "Unpackaged"

# New scratch row added at the bottom
$ws4.Range("A7").Value = 7
$ws4.Range("B7").Value = "'2"
$ws4.Range("C7").Value = "'4"
$ws4.Range("D7").Value = "Unpackaged"
$ws4.Range("E7").Value = 3
$ws4.Range("F7").Value = 5

$ws4.Activate()
$ws4.Range("D6").Select()

# ---------------------------------------------------------------------
# Sheet 5: Dried -> Dried Fruits and Coffee
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "Dried Fruits and Coffee"

$ws5.Range("A2").Value = 1
$ws5.Range("D2").Value = "Packaged"

$ws5.Range("A3").Value = 2
$ws5.Range("D3").Value = "Packaged"

$ws5.Range("A4").Value = 3
$ws5.Range("C4").Value = "bahceden"
$ws5.Range("D4").Value = "Packaged"

$ws5.Range("A5").Value = 4
$ws5.Range("D5").Value = "Packaged"

$ws5.Range("A6").Value = 5
$ws5.Range("D6").Value = "Packaged"

$ws5.Columns.Item(4).ColumnWidth = 11.05
$ws5.Columns.Item(6).ColumnWidth = 12.5

$ws5.Activate()
$ws5.Range("D6").Select()

# ---------------------------------------------------------------------
# Sheet 6: Snacks -> Snack Food
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Name = "Snack Food"

$ws6.Range("A2").Value = 1
$ws6.Range("D2").Value = "Packaged"

$ws6.Range("A3").Value = 2
$ws6.Range("D3").Value = "Packaged"

$ws6.Range("A4").Value = 3
$ws6.Range("D4").Value = "Packaged"

$ws6.Range("A5").Value = 4
$ws6.Range("D5").Value = "Packaged"

$ws6.Range("A6").Value = 5
$ws6.Range("D6").Value = "Packaged"

$ws6.Activate()
$ws6.Rows.Item(7).Select()
